$d = $word.ActiveDocument

# Update the date/title paragraph
$d.Content.Find.Execute("2024-03-14 Thursday", $true, $false, $false, $false, $false, $true, 0, $false, "2024-03-15 Friday", 1)

# Update each arithmetic expression cell in the table, addressed by (row, col)
# so that duplicate expressions map to the correct distinct replacement.
# Use wdReplaceOne (1) and wdFindStop (0) so the replace stays confined to the cell range.
$tbl = $d.Tables.Item(1)

$cell = $tbl.Cell(1, 1)
$cell.Range.Find.Execute("42-39=", $true, $false, $false, $false, $false, $true, 0, $false, "42+39=", 1)
$cell = $tbl.Cell(1, 2)
$cell.Range.Find.Execute("85-76=", $true, $false, $false, $false, $false, $true, 0, $false, "36+49=", 1)
$cell = $tbl.Cell(1, 3)
$cell.Range.Find.Execute("49+32=", $true, $false, $false, $false, $false, $true, 0, $false, "49+4=", 1)
$cell = $tbl.Cell(1, 4)
$cell.Range.Find.Execute("90-7=", $true, $false, $false, $false, $false, $true, 0, $false, "70-38=", 1)
$cell = $tbl.Cell(1, 5)
$cell.Range.Find.Execute("46+38=", $true, $false, $false, $false, $false, $true, 0, $false, "59+33=", 1)
$cell = $tbl.Cell(2, 1)
$cell.Range.Find.Execute("29+66=", $true, $false, $false, $false, $false, $true, 0, $false, "40-12=", 1)
$cell = $tbl.Cell(2, 2)
$cell.Range.Find.Execute("19+15=", $true, $false, $false, $false, $false, $true, 0, $false, "17+54=", 1)
$cell = $tbl.Cell(2, 3)
$cell.Range.Find.Execute("69+3=", $true, $false, $false, $false, $false, $true, 0, $false, "16+78=", 1)
$cell = $tbl.Cell(2, 4)
$cell.Range.Find.Execute("40-9=", $true, $false, $false, $false, $false, $true, 0, $false, "86-67=", 1)
$cell = $tbl.Cell(2, 5)
$cell.Range.Find.Execute("76+6=", $true, $false, $false, $false, $false, $true, 0, $false, "34-27=", 1)
$cell = $tbl.Cell(3, 1)
$cell.Range.Find.Execute("39+57=", $true, $false, $false, $false, $false, $true, 0, $false, "3+78=", 1)
$cell = $tbl.Cell(3, 2)
$cell.Range.Find.Execute("8+88=", $true, $false, $false, $false, $false, $true, 0, $false, "29+59=", 1)
$cell = $tbl.Cell(3, 3)
$cell.Range.Find.Execute("84-68=", $true, $false, $false, $false, $false, $true, 0, $false, "73-14=", 1)
$cell = $tbl.Cell(3, 4)
$cell.Range.Find.Execute("93-88=", $true, $false, $false, $false, $false, $true, 0, $false, "78-39=", 1)
$cell = $tbl.Cell(3, 5)
$cell.Range.Find.Execute("15+19=", $true, $false, $false, $false, $false, $true, 0, $false, "66+18=", 1)
$cell = $tbl.Cell(4, 1)
$cell.Range.Find.Execute("70-9=", $true, $false, $false, $false, $false, $true, 0, $false, "19+68=", 1)
$cell = $tbl.Cell(4, 2)
$cell.Range.Find.Execute("47+47=", $true, $false, $false, $false, $false, $true, 0, $false, "71-16=", 1)
$cell = $tbl.Cell(4, 3)
$cell.Range.Find.Execute("73-18=", $true, $false, $false, $false, $false, $true, 0, $false, "33-17=", 1)
$cell = $tbl.Cell(4, 4)
$cell.Range.Find.Execute("48+38=", $true, $false, $false, $false, $false, $true, 0, $false, "15+36=", 1)
$cell = $tbl.Cell(4, 5)
$cell.Range.Find.Execute("97-49=", $true, $false, $false, $false, $false, $true, 0, $false, "90-1=", 1)
$cell = $tbl.Cell(5, 1)
$cell.Range.Find.Execute("57-8=", $true, $false, $false, $false, $false, $true, 0, $false, "29+22=", 1)
$cell = $tbl.Cell(5, 2)
$cell.Range.Find.Execute("54+27=", $true, $false, $false, $false, $false, $true, 0, $false, "69+29=", 1)
$cell = $tbl.Cell(5, 3)
$cell.Range.Find.Execute("29+24=", $true, $false, $false, $false, $false, $true, 0, $false, "9+75=", 1)
$cell = $tbl.Cell(5, 4)
$cell.Range.Find.Execute("75-16=", $true, $false, $false, $false, $false, $true, 0, $false, "77+5=", 1)
$cell = $tbl.Cell(5, 5)
$cell.Range.Find.Execute("58+13=", $true, $false, $false, $false, $false, $true, 0, $false, "24-19=", 1)
$cell = $tbl.Cell(6, 1)
$cell.Range.Find.Execute("62-56=", $true, $false, $false, $false, $false, $true, 0, $false, "83-17=", 1)
$cell = $tbl.Cell(6, 2)
$cell.Range.Find.Execute("50-29=", $true, $false, $false, $false, $false, $true, 0, $false, "97-78=", 1)
$cell = $tbl.Cell(6, 3)
$cell.Range.Find.Execute("76-39=", $true, $false, $false, $false, $false, $true, 0, $false, "91-65=", 1)
$cell = $tbl.Cell(6, 4)
$cell.Range.Find.Execute("58+37=", $true, $false, $false, $false, $false, $true, 0, $false, "29+62=", 1)
$cell = $tbl.Cell(6, 5)
$cell.Range.Find.Execute("33-9=", $true, $false, $false, $false, $false, $true, 0, $false, "81-6=", 1)
$cell = $tbl.Cell(7, 1)
$cell.Range.Find.Execute("56+9=", $true, $false, $false, $false, $false, $true, 0, $false, "32-25=", 1)
$cell = $tbl.Cell(7, 2)
$cell.Range.Find.Execute("33-27=", $true, $false, $false, $false, $false, $true, 0, $false, "84-48=", 1)
$cell = $tbl.Cell(7, 3)
$cell.Range.Find.Execute("73-68=", $true, $false, $false, $false, $false, $true, 0, $false, "95-57=", 1)
$cell = $tbl.Cell(7, 4)
$cell.Range.Find.Execute("73-17=", $true, $false, $false, $false, $false, $true, 0, $false, "77-49=", 1)
$cell = $tbl.Cell(7, 5)
$cell.Range.Find.Execute("62-34=", $true, $false, $false, $false, $false, $true, 0, $false, "65+8=", 1)
$cell = $tbl.Cell(8, 1)
$cell.Range.Find.Execute("30-22=", $true, $false, $false, $false, $false, $true, 0, $false, "52-6=", 1)
$cell = $tbl.Cell(8, 2)
$cell.Range.Find.Execute("86-27=", $true, $false, $false, $false, $false, $true, 0, $false, "66-57=", 1)
$cell = $tbl.Cell(8, 3)
$cell.Range.Find.Execute("82-43=", $true, $false, $false, $false, $false, $true, 0, $false, "74-55=", 1)
$cell = $tbl.Cell(8, 4)
$cell.Range.Find.Execute("57-9=", $true, $false, $false, $false, $false, $true, 0, $false, "12-4=", 1)
$cell = $tbl.Cell(8, 5)
$cell.Range.Find.Execute("29+16=", $true, $false, $false, $false, $false, $true, 0, $false, "17+39=", 1)
$cell = $tbl.Cell(9, 1)
$cell.Range.Find.Execute("5+78=", $true, $false, $false, $false, $false, $true, 0, $false, "95-9=", 1)
$cell = $tbl.Cell(9, 2)
$cell.Range.Find.Execute("87+4=", $true, $false, $false, $false, $false, $true, 0, $false, "67+25=", 1)
$cell = $tbl.Cell(9, 3)
$cell.Range.Find.Execute("28+46=", $true, $false, $false, $false, $false, $true, 0, $false, "80-31=", 1)
$cell = $tbl.Cell(9, 4)
$cell.Range.Find.Execute("64-48=", $true, $false, $false, $false, $false, $true, 0, $false, "45+19=", 1)
$cell = $tbl.Cell(9, 5)
$cell.Range.Find.Execute("92-68=", $true, $false, $false, $false, $false, $true, 0, $false, "77-59=", 1)
$cell = $tbl.Cell(10, 1)
$cell.Range.Find.Execute("56-48=", $true, $false, $false, $false, $false, $true, 0, $false, "6+55=", 1)
$cell = $tbl.Cell(10, 2)
$cell.Range.Find.Execute("92-33=", $true, $false, $false, $false, $false, $true, 0, $false, "35-29=", 1)
$cell = $tbl.Cell(10, 3)
$cell.Range.Find.Execute("73-34=", $true, $false, $false, $false, $false, $true, 0, $false, "36+49=", 1)
$cell = $tbl.Cell(10, 4)
$cell.Range.Find.Execute("19+38=", $true, $false, $false, $false, $false, $true, 0, $false, "9+58=", 1)
$cell = $tbl.Cell(10, 5)
$cell.Range.Find.Execute("89+8=", $true, $false, $false, $false, $false, $true, 0, $false, "48+18=", 1)
$cell = $tbl.Cell(11, 1)
$cell.Range.Find.Execute("93-7=", $true, $false, $false, $false, $false, $true, 0, $false, "50-45=", 1)
$cell = $tbl.Cell(11, 2)
$cell.Range.Find.Execute("4+19=", $true, $false, $false, $false, $false, $true, 0, $false, "94-27=", 1)
$cell = $tbl.Cell(11, 3)
$cell.Range.Find.Execute("16+45=", $true, $false, $false, $false, $false, $true, 0, $false, "43-14=", 1)
$cell = $tbl.Cell(11, 4)
$cell.Range.Find.Execute("80-78=", $true, $false, $false, $false, $false, $true, 0, $false, "93-64=", 1)
$cell = $tbl.Cell(11, 5)
$cell.Range.Find.Execute("64+27=", $true, $false, $false, $false, $false, $true, 0, $false, "19+54=", 1)
$cell = $tbl.Cell(12, 1)
$cell.Range.Find.Execute("91-22=", $true, $false, $false, $false, $false, $true, 0, $false, "88-19=", 1)
$cell = $tbl.Cell(12, 2)
$cell.Range.Find.Execute("23+8=", $true, $false, $false, $false, $false, $true, 0, $false, "70-31=", 1)
$cell = $tbl.Cell(12, 3)
$cell.Range.Find.Execute("72-17=", $true, $false, $false, $false, $false, $true, 0, $false, "53-38=", 1)
$cell = $tbl.Cell(12, 4)
$cell.Range.Find.Execute("51-23=", $true, $false, $false, $false, $false, $true, 0, $false, "73-65=", 1)
$cell = $tbl.Cell(12, 5)
$cell.Range.Find.Execute("91-66=", $true, $false, $false, $false, $false, $true, 0, $false, "46+6=", 1)
$cell = $tbl.Cell(13, 1)
$cell.Range.Find.Execute("59+37=", $true, $false, $false, $false, $false, $true, 0, $false, "88+5=", 1)
$cell = $tbl.Cell(13, 2)
$cell.Range.Find.Execute("53+8=", $true, $false, $false, $false, $false, $true, 0, $false, "91-79=", 1)
$cell = $tbl.Cell(13, 3)
$cell.Range.Find.Execute("7+67=", $true, $false, $false, $false, $false, $true, 0, $false, "57-48=", 1)
$cell = $tbl.Cell(13, 4)
$cell.Range.Find.Execute("52-16=", $true, $false, $false, $false, $false, $true, 0, $false, "91-16=", 1)
$cell = $tbl.Cell(13, 5)
$cell.Range.Find.Execute("72-64=", $true, $false, $false, $false, $false, $true, 0, $false, "80-3=", 1)
$cell = $tbl.Cell(14, 1)
$cell.Range.Find.Execute("8+6=", $true, $false, $false, $false, $false, $true, 0, $false, "81-78=", 1)
$cell = $tbl.Cell(14, 2)
$cell.Range.Find.Execute("6+28=", $true, $false, $false, $false, $false, $true, 0, $false, "18+68=", 1)
$cell = $tbl.Cell(14, 3)
$cell.Range.Find.Execute("83-29=", $true, $false, $false, $false, $false, $true, 0, $false, "25+37=", 1)
$cell = $tbl.Cell(14, 4)
$cell.Range.Find.Execute("3+48=", $true, $false, $false, $false, $false, $true, 0, $false, "60-32=", 1)
$cell = $tbl.Cell(14, 5)
$cell.Range.Find.Execute("82-43=", $true, $false, $false, $false, $false, $true, 0, $false, "46+26=", 1)
$cell = $tbl.Cell(15, 1)
$cell.Range.Find.Execute("22-4=", $true, $false, $false, $false, $false, $true, 0, $false, "9+13=", 1)
$cell = $tbl.Cell(15, 2)
$cell.Range.Find.Execute("19+17=", $true, $false, $false, $false, $false, $true, 0, $false, "68+29=", 1)
$cell = $tbl.Cell(15, 3)
$cell.Range.Find.Execute("83-58=", $true, $false, $false, $false, $false, $true, 0, $false, "33-16=", 1)
$cell = $tbl.Cell(15, 4)
$cell.Range.Find.Execute("44-6=", $true, $false, $false, $false, $false, $true, 0, $false, "83-67=", 1)
$cell = $tbl.Cell(15, 5)
$cell.Range.Find.Execute("26+57=", $true, $false, $false, $false, $false, $true, 0, $false, "52-46=", 1)
$cell = $tbl.Cell(16, 1)
$cell.Range.Find.Execute("92-35=", $true, $false, $false, $false, $false, $true, 0, $false, "50-38=", 1)
$cell = $tbl.Cell(16, 2)
$cell.Range.Find.Execute("67+14=", $true, $false, $false, $false, $false, $true, 0, $false, "17+15=", 1)
$cell = $tbl.Cell(16, 3)
$cell.Range.Find.Execute("72-15=", $true, $false, $false, $false, $false, $true, 0, $false, "90-14=", 1)
$cell = $tbl.Cell(16, 4)
$cell.Range.Find.Execute("7+38=", $true, $false, $false, $false, $false, $true, 0, $false, "90-88=", 1)
$cell = $tbl.Cell(16, 5)
$cell.Range.Find.Execute("71-49=", $true, $false, $false, $false, $false, $true, 0, $false, "78+18=", 1)
$cell = $tbl.Cell(17, 1)
$cell.Range.Find.Execute("34+29=", $true, $false, $false, $false, $false, $true, 0, $false, "38+28=", 1)
$cell = $tbl.Cell(17, 2)
$cell.Range.Find.Execute("92-58=", $true, $false, $false, $false, $false, $true, 0, $false, "12+59=", 1)
$cell = $tbl.Cell(17, 3)
$cell.Range.Find.Execute("19+48=", $true, $false, $false, $false, $false, $true, 0, $false, "19+75=", 1)
$cell = $tbl.Cell(17, 4)
$cell.Range.Find.Execute("47+6=", $true, $false, $false, $false, $false, $true, 0, $false, "69+13=", 1)
$cell = $tbl.Cell(17, 5)
$cell.Range.Find.Execute("6+9=", $true, $false, $false, $false, $false, $true, 0, $false, "2+79=", 1)
$cell = $tbl.Cell(18, 1)
$cell.Range.Find.Execute("18+48=", $true, $false, $false, $false, $false, $true, 0, $false, "39+35=", 1)
$cell = $tbl.Cell(18, 2)
$cell.Range.Find.Execute("90-53=", $true, $false, $false, $false, $false, $true, 0, $false, "39+16=", 1)
$cell = $tbl.Cell(18, 3)
$cell.Range.Find.Execute("13-9=", $true, $false, $false, $false, $false, $true, 0, $false, "14+48=", 1)
$cell = $tbl.Cell(18, 4)
$cell.Range.Find.Execute("36+57=", $true, $false, $false, $false, $false, $true, 0, $false, "37-8=", 1)
$cell = $tbl.Cell(18, 5)
$cell.Range.Find.Execute("65+17=", $true, $false, $false, $false, $false, $true, 0, $false, "27+59=", 1)
$cell = $tbl.Cell(19, 1)
$cell.Range.Find.Execute("29+63=", $true, $false, $false, $false, $false, $true, 0, $false, "52-43=", 1)
$cell = $tbl.Cell(19, 2)
$cell.Range.Find.Execute("58+28=", $true, $false, $false, $false, $false, $true, 0, $false, "61-52=", 1)
$cell = $tbl.Cell(19, 3)
$cell.Range.Find.Execute("55+8=", $true, $false, $false, $false, $false, $true, 0, $false, "8+86=", 1)
$cell = $tbl.Cell(19, 4)
$cell.Range.Find.Execute("8+14=", $true, $false, $false, $false, $false, $true, 0, $false, "72-13=", 1)
$cell = $tbl.Cell(19, 5)
$cell.Range.Find.Execute("64-57=", $true, $false, $false, $false, $false, $true, 0, $false, "9+18=", 1)
$cell = $tbl.Cell(20, 1)
$cell.Range.Find.Execute("44-19=", $true, $false, $false, $false, $false, $true, 0, $false, "45+8=", 1)
$cell = $tbl.Cell(20, 2)
$cell.Range.Find.Execute("70-9=", $true, $false, $false, $false, $false, $true, 0, $false, "90-35=", 1)
$cell = $tbl.Cell(20, 3)
$cell.Range.Find.Execute("41-37=", $true, $false, $false, $false, $false, $true, 0, $false, "44-29=", 1)
$cell = $tbl.Cell(20, 4)
$cell.Range.Find.Execute("95-38=", $true, $false, $false, $false, $false, $true, 0, $false, "92-34=", 1)
$cell = $tbl.Cell(20, 5)
$cell.Range.Find.Execute("30-21=", $true, $false, $false, $false, $false, $true, 0, $false, "70-18=", 1)
